# Refresh the cryptos price/volume table (columns D "Price" and E "Volume(1h)")
# with the latest scraped values. Price/Volume cells are stored as plain text
# (not numbers) in this sheet, so for D-column values that look like numbers
# we briefly force Text format before assigning, then clear the format again
# so the cell's style matches the rest of the untouched column (no stray
# NumberFormat artifacts left behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.716.41"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "3.714.82"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "672.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.14"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.445"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.98%  "

$ws.Range("D14").Value = "3.743.65"
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("D15").Value = "69.733.32"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("E16").Value = "  +1.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "474.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.656"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").Value = "3.861.28"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000129"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.39%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.01"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.22%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.96"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("D35").Value = "3.703.63"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.57"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.70%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.26"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0913"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.38"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.51%  "

$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000285"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.07%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("E51").Value = "  +1.63%  "
